# ---------------------------------------------------------------------------
# Applies the "Notes" section rewrite described by the commit:
#   "Seedboxes now of bird child object when carried Tag whole house with
#    CanUseClaw"
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Get-ParaIndexByExactText($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# 1) "Notes" heading: add a new "Can hide food In bathroom bin" paragraph right
#    after it. Re-stamping the "Notes" run's text (identical content) drops the
#    stale <w:lastRenderedPageBreak/> marker, same as Word does when it next
#    lays out the page and re-serialises that run.
$headingRng = $d.Content
$headingRng.Find.Execute("Notes", $true, $false, $false, $false, $false, $true, 1, $false, "Notes", 2) | Out-Null

$idx = Get-ParaIndexByExactText $d "Notes"
$p = $d.Paragraphs($idx)
$p.Range.InsertParagraphAfter()
$newP = $d.Paragraphs($idx + 1)
$newP.Style = "Normal"
$newP.Range.Text = "Can hide food In bathroom bin"

# 2) Remove "Turn baked lighting back on for builds (Baked global Illumination)"
#    from just below TODO (it reappears later, near "Cage needs colliders").
$idx = Get-ParaIndexByExactText $d "Turn baked lighting back on for builds (Baked global Illumination)"
$p = $d.Paragraphs($idx)
$d.Range($p.Range.Start, $p.Range.End).Delete()

# 3) Remove the "4th Milestone" paragraph that used to sit below TODO.
$idx = Get-ParaIndexByExactText $d "4th Milestone"
$p = $d.Paragraphs($idx)
$d.Range($p.Range.Start, $p.Range.End).Delete()

# 4) Add "Sometimes when letting go of seed, physics goes a bit mad." after
#    "Disable player control (including partial) when eating/bathing/in tutorial".
$idx = Get-ParaIndexByExactText $d "Disable player control (including partial) when eating/bathing/in tutorial"
$p = $d.Paragraphs($idx)
$p.Range.InsertParagraphAfter()
$newP = $d.Paragraphs($idx + 1)
$newP.Style = "Normal"
$newP.Range.Text = "Sometimes when letting go of seed, physics goes a bit mad."

# 5) DOING section: "3rd Milestone" -> "4th Milestone"
$renameRng = $d.Content
$renameRng.Find.Execute("3rd Milestone", $true, $false, $false, $false, $false, $true, 1, $false, "4th Milestone", 2) | Out-Null

# 6) "PlayerController"/" now checks if " runs (with surrounding proofErr spell
#    markers) must collapse into a single run, while the following "collision"
#    run (and everything after it) keeps its own formatting untouched. A plain
#    Find/Replace merges every run from the edit point through the end of the
#    paragraph, so the paragraph is first split right before "collision",
#    the now-isolated lead-in sentence is rewritten in isolation, and then the
#    split is rejoined.
$splitRng = $d.Content
$splitRng.Find.Execute("collision", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $splitRng.Start
$d.Range($splitPos, $splitPos).InsertParagraphBefore()

$idx = Get-ParaIndexByExactText $d "PlayerController now checks if "
$p = $d.Paragraphs($idx)
$d.Range($p.Range.Start, $p.Range.End).Delete()
$d.Paragraphs($idx - 1).Range.InsertParagraphAfter()
$newP = $d.Paragraphs($idx)
$newP.Style = "Normal"
$newP.Range.Text = "PlayerController now checks if "

# Rejoin: delete the paragraph mark that now separates the lead-in sentence
# from the "collision ..." remainder.
$joinPos = $d.Paragraphs($idx).Range.End - 1
$d.Range($joinPos, $joinPos + 1).Delete()

# 7) "Load/" + "SaveProgress" runs (with proofErr spell markers) -> one run.
$idx = Get-ParaIndexByExactText $d "Load/SaveProgress"
$p = $d.Paragraphs($idx)
$d.Range($p.Range.Start, $p.Range.End).Delete()
$d.Paragraphs($idx - 1).Range.InsertParagraphAfter()
$newP = $d.Paragraphs($idx)
$newP.Style = "Normal"
$newP.Range.Text = "Load/SaveProgress"

# 8) Append "Turn baked lighting back on for builds (Baked global Illumination)"
#    and "3rd Milestone" after "Cage needs colliders", followed by one more
#    blank paragraph (on top of the blank paragraph already there).
$idx = Get-ParaIndexByExactText $d "Cage needs colliders"
$p = $d.Paragraphs($idx)
$p.Range.InsertParagraphAfter()
$newP = $d.Paragraphs($idx + 1)
$newP.Style = "Normal"
$newP.Range.Text = "Turn baked lighting back on for builds (Baked global Illumination)"

$p2 = $d.Paragraphs($idx + 1)
$p2.Range.InsertParagraphAfter()
$newP2 = $d.Paragraphs($idx + 2)
$newP2.Style = "Normal"
$newP2.Range.Text = "3rd Milestone"

$p3 = $d.Paragraphs($idx + 2)
$p3.Range.InsertParagraphAfter()
$newP3 = $d.Paragraphs($idx + 3)
$newP3.Style = "Normal"

Write-Host "done phase 1"
